# Incorporación de Footy Stats DB
#
# Adds a second "moneyline" (H/D/A) odds block to Sheet2 (columns J:L,
# plus the margin helper in N3:N4), and converts the already-existing
# per-row "=1/Cxx" formulas in Atlas_Tigres_ida!D30:D58 into one shared
# formula. Also updates sheet selection / active-tab state to match
# where the author ended up (Sheet2 active, selection on L6; Atlas_Tigres_ida
# selection moved to H1).

$wb = $excel.ActiveWorkbook
$wsOdds = $wb.Worksheets.Item("Atlas_Tigres_ida")
$wsFooty = $wb.Worksheets.Item("Sheet2")

# --- Atlas_Tigres_ida: turn the D30:D58 per-cell formulas into a shared formula ---
$wsOdds.Range("D30:D58").Formula = "=1/C30"

# Move the (inactive) selection on Atlas_Tigres_ida to H1.
$wsOdds.Range("H1").Select()

# --- Sheet2: new Footy Stats DB block (H/D/A moneyline odds) ---
$wsFooty.Range("J2").Value = "H"
$wsFooty.Range("K2").Value = "D"
$wsFooty.Range("L2").Value = "A"

$wsFooty.Range("J3").Value = 1.5
$wsFooty.Range("K3").Value = 5
$wsFooty.Range("L3").Value = 7.5
$wsFooty.Range("N3").Value = "Margin"

$wsFooty.Range("J4").Formula = "=1/J3"
$wsFooty.Range("K4:L4").Formula = "=1/K3"
$wsFooty.Range("N4").Value = 0.08

$wsFooty.Range("J5").Formula = "=J4+`$N`$4/3"
$wsFooty.Range("K5:L5").Formula = "=K4+`$N`$4/3"

$wsFooty.Range("J6").Formula = "=1/J5"
$wsFooty.Range("K6:L6").Formula = "=1/K5"

# N4 carries the 3-decimal margin format (numFmt 0.000) - apply last so it
# doesn't leak onto the cells created afterwards.
$wsFooty.Range("N4").NumberFormat = "0.000"

# Sheet2 becomes the active sheet/tab, with the cursor left on L6.
$wsFooty.Activate()
$wsFooty.Range("L6").Select()
